# Generate Report for Handoff
#
# The HO xliff generation produced a new Xliff/markdown GUID and pushed a
# newer "Latest HO Xliff Generate Date" / handoff timestamp to the report.
# Update every cell that carries the old GUID / old timestamps (and the
# matching hyperlink display text) on all three report sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "f1f48d12-dc06-4ec2-902b-72366931228f"
$newGuid = "e52db018-6830-4abf-ba72-dd23dfc01521"

$oldHash = "a26b4cc3da0d181a7984a25ad42f1a9cc2bc0471"
$newHash = "c7f9f1f55fc52c0b745fea7ef0e3f7206b93ab5a"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-21 03:01:02"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-21 03:00:56"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
